$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 386 (shift old rows 386-460 down to 388-462)
$ws.Rows.Item(386).Resize(2).Insert()

# Row 384: new data point (previously held "1a nueva(o)" O'Higgins data dated 44505)
$ws.Cells.Item(384, 4).Value2 = 44694
$ws.Cells.Item(384, 9).Value2 = "1a (cosecha)"
$ws.Cells.Item(384, 10).Value2 = 600
$ws.Cells.Item(384, 11).Value2 = 6000
$ws.Cells.Item(384, 12).Value2 = 6500
$ws.Cells.Item(384, 13).Value2 = 6250
$ws.Cells.Item(384, 16).Value2 = 347

# Row 385: new data point (previously held "2a nueva(o)" O'Higgins data dated 44505)
$ws.Cells.Item(385, 4).Value2 = 44694
$ws.Cells.Item(385, 9).Value2 = "2a (cosecha)"
$ws.Cells.Item(385, 10).Value2 = 300
$ws.Cells.Item(385, 11).Value2 = 5500
$ws.Cells.Item(385, 12).Value2 = 5500
$ws.Cells.Item(385, 13).Value2 = 5500
$ws.Cells.Item(385, 16).Value2 = 306

# Row 386: fill newly-inserted blank row with data that used to be row 384
$ws.Cells.Item(386, 1).Value2 = 11
$ws.Cells.Item(386, 2).Value2 = "Vega Monumental Concepción"
$ws.Cells.Item(386, 3).Value2 = "Bíobío"
$ws.Cells.Item(386, 4).Value2 = 44505
$ws.Cells.Item(386, 5).Value2 = 8
$ws.Cells.Item(386, 6).Value2 = 100112004
$ws.Cells.Item(386, 7).Value2 = "Cebolla"
$ws.Cells.Item(386, 8).Value2 = "Sin especificar"
$ws.Cells.Item(386, 9).Value2 = "1a nueva(o)"
$ws.Cells.Item(386, 10).Value2 = 430
$ws.Cells.Item(386, 11).Value2 = 5000
$ws.Cells.Item(386, 12).Value2 = 5500
$ws.Cells.Item(386, 13).Value2 = 5291
$ws.Cells.Item(386, 14).Value2 = "$/malla 18 kilos"
$ws.Cells.Item(386, 15).Value2 = "Región de O'Higgins"
$ws.Cells.Item(386, 16).Value2 = 294
$ws.Cells.Item(386, 17).Value2 = 18
$ws.Cells.Item(386, 18).Value2 = "Hortaliza"

# Row 387: fill newly-inserted blank row with data that used to be row 385
$ws.Cells.Item(387, 1).Value2 = 11
$ws.Cells.Item(387, 2).Value2 = "Vega Monumental Concepción"
$ws.Cells.Item(387, 3).Value2 = "Bíobío"
$ws.Cells.Item(387, 4).Value2 = 44505
$ws.Cells.Item(387, 5).Value2 = 8
$ws.Cells.Item(387, 6).Value2 = 100112004
$ws.Cells.Item(387, 7).Value2 = "Cebolla"
$ws.Cells.Item(387, 8).Value2 = "Sin especificar"
$ws.Cells.Item(387, 9).Value2 = "2a nueva(o)"
$ws.Cells.Item(387, 10).Value2 = 200
$ws.Cells.Item(387, 11).Value2 = 4500
$ws.Cells.Item(387, 12).Value2 = 4500
$ws.Cells.Item(387, 13).Value2 = 4500
$ws.Cells.Item(387, 14).Value2 = "$/malla 18 kilos"
$ws.Cells.Item(387, 15).Value2 = "Región de O'Higgins"
$ws.Cells.Item(387, 16).Value2 = 250
$ws.Cells.Item(387, 17).Value2 = 18
$ws.Cells.Item(387, 18).Value2 = "Hortaliza"
